$wb = $excel.ActiveWorkbook
$models = $wb.Worksheets.Item("Models")
$instances = $wb.Worksheets.Item("Instances")

# The "Path" column was actually the parent path, not a full path - rename the
# header and make the "Level" cell pull the tag name live from the Models sheet.
$instances.Range("C3").Value = "Parent"
$instances.Range("E3").Formula = "=Models!C3"

# Tidy up capitalization of the Clsid header and drop the curly braces Excel
# doesn't need around the literal GUID value in the example row.
$instances.Range("K3").Value = "ClsId"
$instances.Range("K7").Value = "91210ec1-58ac-41f9-b840-b39b965076fc"

# Restore the selections/active cells to match the saved view.
$models.Range("D2").Select()
$instances.Activate()
$instances.Range("E3").Select()
